$wb = $excel.ActiveWorkbook

# --- Insert 'Rebounds' sheet right after 'Assists' ---
$afterAssists = $wb.Worksheets.Item("Assists")
$wsRebounds = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterAssists)
$wsRebounds.Name = "Rebounds"

$reboundsData = New-Object 'object[,]' 11,17
$reboundsData[0,0] = "Game Time (PST)"
$reboundsData[0,1] = "Opponent"
$reboundsData[0,2] = "Taylor Hendricks"
$reboundsData[0,3] = "Kyle Anderson"
$reboundsData[0,4] = "Keyonte George"
$reboundsData[0,5] = "Cody Williams"
$reboundsData[0,6] = "Isaiah Collier"
$reboundsData[0,7] = "Svi Mykhailiuk"
$reboundsData[0,8] = "Walter Clayton Jr."
$reboundsData[0,9] = "Elijah Harkless"
$reboundsData[0,10] = "Ace Bailey"
$reboundsData[0,11] = "Kyle Filipowski"
$reboundsData[0,12] = "Lauri Markkanen"
$reboundsData[0,13] = "Walker Kessler"
$reboundsData[0,14] = "Brice Sensabaugh"
$reboundsData[0,15] = "Jusuf Nurkić"
$reboundsData[0,16] = "Kevin Love"
$reboundsData[1,0] = "2025-10-22"
$reboundsData[1,1] = "LAC"
$reboundsData[1,2] = 6
$reboundsData[1,3] = 0
$reboundsData[1,4] = 2
$reboundsData[1,5] = 0
$reboundsData[1,6] = 0
$reboundsData[1,7] = 1
$reboundsData[1,8] = 6
$reboundsData[1,9] = 0
$reboundsData[1,10] = 4
$reboundsData[1,11] = 4
$reboundsData[1,12] = 6
$reboundsData[1,13] = 9
$reboundsData[1,14] = 4
$reboundsData[1,15] = 4
$reboundsData[1,16] = 0
$reboundsData[2,0] = "2025-10-24"
$reboundsData[2,1] = "SAC"
$reboundsData[2,2] = 3
$reboundsData[2,3] = 0
$reboundsData[2,4] = 2
$reboundsData[2,5] = 0
$reboundsData[2,6] = 0
$reboundsData[2,7] = 4
$reboundsData[2,8] = 5
$reboundsData[2,9] = 0
$reboundsData[2,10] = 2
$reboundsData[2,11] = 6
$reboundsData[2,12] = 4
$reboundsData[2,13] = 9
$reboundsData[2,14] = 4
$reboundsData[2,15] = 11
$reboundsData[2,16] = 0
$reboundsData[3,0] = "2025-10-27"
$reboundsData[3,1] = "PHX"
$reboundsData[3,2] = 6
$reboundsData[3,3] = 0
$reboundsData[3,4] = 5
$reboundsData[3,5] = 0
$reboundsData[3,6] = 0
$reboundsData[3,7] = 3
$reboundsData[3,8] = 3
$reboundsData[3,9] = 0
$reboundsData[3,10] = 1
$reboundsData[3,11] = 6
$reboundsData[3,12] = 14
$reboundsData[3,13] = 11
$reboundsData[3,14] = 2
$reboundsData[3,15] = 13
$reboundsData[3,16] = 0
$reboundsData[4,0] = "2025-10-29"
$reboundsData[4,1] = "POR"
$reboundsData[4,2] = 4
$reboundsData[4,3] = 0
$reboundsData[4,4] = 2
$reboundsData[4,5] = 1
$reboundsData[4,6] = 0
$reboundsData[4,7] = 2
$reboundsData[4,8] = 1
$reboundsData[4,9] = 3
$reboundsData[4,10] = 3
$reboundsData[4,11] = 4
$reboundsData[4,12] = 5
$reboundsData[4,13] = 12
$reboundsData[4,14] = 5
$reboundsData[4,15] = 4
$reboundsData[4,16] = 0
$reboundsData[5,0] = "2025-10-31"
$reboundsData[5,1] = "PHX"
$reboundsData[5,2] = 0
$reboundsData[5,3] = 3
$reboundsData[5,4] = 8
$reboundsData[5,5] = 1
$reboundsData[5,6] = 0
$reboundsData[5,7] = 3
$reboundsData[5,8] = 0
$reboundsData[5,9] = 0
$reboundsData[5,10] = 8
$reboundsData[5,11] = 2
$reboundsData[5,12] = 4
$reboundsData[5,13] = 13
$reboundsData[5,14] = 0
$reboundsData[5,15] = 6
$reboundsData[5,16] = 0
$reboundsData[6,0] = "2025-11-02"
$reboundsData[6,1] = "CHA"
$reboundsData[6,2] = 7
$reboundsData[6,3] = 0
$reboundsData[6,4] = 2
$reboundsData[6,5] = 2
$reboundsData[6,6] = 0
$reboundsData[6,7] = 2
$reboundsData[6,8] = 2
$reboundsData[6,9] = 0
$reboundsData[6,10] = 0
$reboundsData[6,11] = 1
$reboundsData[6,12] = 7
$reboundsData[6,13] = 0
$reboundsData[6,14] = 1
$reboundsData[6,15] = 9
$reboundsData[6,16] = 5
$reboundsData[7,0] = "2025-11-03"
$reboundsData[7,1] = "BOS"
$reboundsData[7,2] = 3
$reboundsData[7,3] = 0
$reboundsData[7,4] = 5
$reboundsData[7,5] = 0
$reboundsData[7,6] = 0
$reboundsData[7,7] = 6
$reboundsData[7,8] = 3
$reboundsData[7,9] = 1
$reboundsData[7,10] = 2
$reboundsData[7,11] = 8
$reboundsData[7,12] = 9
$reboundsData[7,13] = 0
$reboundsData[7,14] = 0
$reboundsData[7,15] = 11
$reboundsData[7,16] = 7
$reboundsData[8,0] = "2025-11-05"
$reboundsData[8,1] = "DET"
$reboundsData[8,2] = 7
$reboundsData[8,3] = 0
$reboundsData[8,4] = 2
$reboundsData[8,5] = 0
$reboundsData[8,6] = 0
$reboundsData[8,7] = 2
$reboundsData[8,8] = 0
$reboundsData[8,9] = 1
$reboundsData[8,10] = 4
$reboundsData[8,11] = 4
$reboundsData[8,12] = 4
$reboundsData[8,13] = 0
$reboundsData[8,14] = 0
$reboundsData[8,15] = 17
$reboundsData[8,16] = 3
$reboundsData[9,0] = "2025-11-07"
$reboundsData[9,1] = "MIN"
$reboundsData[9,2] = 2
$reboundsData[9,3] = 4
$reboundsData[9,4] = 4
$reboundsData[9,5] = 1
$reboundsData[9,6] = 0
$reboundsData[9,7] = 2
$reboundsData[9,8] = 3
$reboundsData[9,9] = 1
$reboundsData[9,10] = 6
$reboundsData[9,11] = 10
$reboundsData[9,12] = 2
$reboundsData[9,13] = 0
$reboundsData[9,14] = 0
$reboundsData[9,15] = 4
$reboundsData[9,16] = 1
$reboundsData[10,0] = "2025-11-10"
$reboundsData[10,1] = "MIN"
$reboundsData[10,2] = 0
$reboundsData[10,3] = 0
$reboundsData[10,4] = 6
$reboundsData[10,5] = 0
$reboundsData[10,6] = 3
$reboundsData[10,7] = 3
$reboundsData[10,8] = 3
$reboundsData[10,9] = 0
$reboundsData[10,10] = 7
$reboundsData[10,11] = 3
$reboundsData[10,12] = 8
$reboundsData[10,13] = 0
$reboundsData[10,14] = 4
$reboundsData[10,15] = 10
$reboundsData[10,16] = 3
$wsRebounds.Range("A1:A11").NumberFormat = "@"
$wsRebounds.Range("A1:Q11").Value = $reboundsData
$wsRebounds.Range("A1:Q1").Font.Bold = $true
$wsRebounds.Range("A1:Q1").HorizontalAlignment = -4108
$wsRebounds.Range("A1:Q1").VerticalAlignment = -4160
$wsRebounds.Range("A1:Q1").Borders.LineStyle = 1

# --- Insert '3PM' sheet right after 'Rebounds' ---
$ws3pm = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsRebounds)
$ws3pm.Name = "3PM"

$threePmData = New-Object 'object[,]' 11,17
$threePmData[0,0] = "Game Time (PST)"
$threePmData[0,1] = "Opponent"
$threePmData[0,2] = "Taylor Hendricks"
$threePmData[0,3] = "Kyle Anderson"
$threePmData[0,4] = "Keyonte George"
$threePmData[0,5] = "Cody Williams"
$threePmData[0,6] = "Isaiah Collier"
$threePmData[0,7] = "Svi Mykhailiuk"
$threePmData[0,8] = "Walter Clayton Jr."
$threePmData[0,9] = "Elijah Harkless"
$threePmData[0,10] = "Ace Bailey"
$threePmData[0,11] = "Kyle Filipowski"
$threePmData[0,12] = "Lauri Markkanen"
$threePmData[0,13] = "Walker Kessler"
$threePmData[0,14] = "Brice Sensabaugh"
$threePmData[0,15] = "Jusuf Nurkić"
$threePmData[0,16] = "Kevin Love"
$threePmData[1,0] = "2025-10-22"
$threePmData[1,1] = "LAC"
$threePmData[1,2] = 2
$threePmData[1,3] = 0
$threePmData[1,4] = 1
$threePmData[1,5] = 0
$threePmData[1,6] = 0
$threePmData[1,7] = 3
$threePmData[1,8] = 1
$threePmData[1,9] = 0
$threePmData[1,10] = 0
$threePmData[1,11] = 1
$threePmData[1,12] = 4
$threePmData[1,13] = 2
$threePmData[1,14] = 3
$threePmData[1,15] = 0
$threePmData[1,16] = 0
$threePmData[2,0] = "2025-10-24"
$threePmData[2,1] = "SAC"
$threePmData[2,2] = 0
$threePmData[2,3] = 0
$threePmData[2,4] = 1
$threePmData[2,5] = 0
$threePmData[2,6] = 0
$threePmData[2,7] = 2
$threePmData[2,8] = 1
$threePmData[2,9] = 0
$threePmData[2,10] = 0
$threePmData[2,11] = 2
$threePmData[2,12] = 4
$threePmData[2,13] = 0
$threePmData[2,14] = 3
$threePmData[2,15] = 0
$threePmData[2,16] = 0
$threePmData[3,0] = "2025-10-27"
$threePmData[3,1] = "PHX"
$threePmData[3,2] = 0
$threePmData[3,3] = 0
$threePmData[3,4] = 2
$threePmData[3,5] = 0
$threePmData[3,6] = 0
$threePmData[3,7] = 2
$threePmData[3,8] = 2
$threePmData[3,9] = 0
$threePmData[3,10] = 2
$threePmData[3,11] = 1
$threePmData[3,12] = 6
$threePmData[3,13] = 1
$threePmData[3,14] = 0
$threePmData[3,15] = 0
$threePmData[3,16] = 0
$threePmData[4,0] = "2025-10-29"
$threePmData[4,1] = "POR"
$threePmData[4,2] = 0
$threePmData[4,3] = 0
$threePmData[4,4] = 2
$threePmData[4,5] = 0
$threePmData[4,6] = 0
$threePmData[4,7] = 1
$threePmData[4,8] = 0
$threePmData[4,9] = 0
$threePmData[4,10] = 0
$threePmData[4,11] = 1
$threePmData[4,12] = 2
$threePmData[4,13] = 3
$threePmData[4,14] = 1
$threePmData[4,15] = 1
$threePmData[4,16] = 0
$threePmData[5,0] = "2025-10-31"
$threePmData[5,1] = "PHX"
$threePmData[5,2] = 0
$threePmData[5,3] = 0
$threePmData[5,4] = 2
$threePmData[5,5] = 2
$threePmData[5,6] = 0
$threePmData[5,7] = 0
$threePmData[5,8] = 1
$threePmData[5,9] = 1
$threePmData[5,10] = 0
$threePmData[5,11] = 0
$threePmData[5,12] = 4
$threePmData[5,13] = 0
$threePmData[5,14] = 1
$threePmData[5,15] = 1
$threePmData[5,16] = 0
$threePmData[6,0] = "2025-11-02"
$threePmData[6,1] = "CHA"
$threePmData[6,2] = 2
$threePmData[6,3] = 0
$threePmData[6,4] = 1
$threePmData[6,5] = 0
$threePmData[6,6] = 0
$threePmData[6,7] = 1
$threePmData[6,8] = 0
$threePmData[6,9] = 0
$threePmData[6,10] = 0
$threePmData[6,11] = 0
$threePmData[6,12] = 3
$threePmData[6,13] = 0
$threePmData[6,14] = 0
$threePmData[6,15] = 0
$threePmData[6,16] = 0
$threePmData[7,0] = "2025-11-03"
$threePmData[7,1] = "BOS"
$threePmData[7,2] = 0
$threePmData[7,3] = 0
$threePmData[7,4] = 2
$threePmData[7,5] = 0
$threePmData[7,6] = 0
$threePmData[7,7] = 1
$threePmData[7,8] = 1
$threePmData[7,9] = 0
$threePmData[7,10] = 0
$threePmData[7,11] = 2
$threePmData[7,12] = 2
$threePmData[7,13] = 0
$threePmData[7,14] = 0
$threePmData[7,15] = 0
$threePmData[7,16] = 1
$threePmData[8,0] = "2025-11-05"
$threePmData[8,1] = "DET"
$threePmData[8,2] = 0
$threePmData[8,3] = 0
$threePmData[8,4] = 1
$threePmData[8,5] = 0
$threePmData[8,6] = 0
$threePmData[8,7] = 6
$threePmData[8,8] = 1
$threePmData[8,9] = 1
$threePmData[8,10] = 0
$threePmData[8,11] = 0
$threePmData[8,12] = 3
$threePmData[8,13] = 0
$threePmData[8,14] = 0
$threePmData[8,15] = 0
$threePmData[8,16] = 0
$threePmData[9,0] = "2025-11-07"
$threePmData[9,1] = "MIN"
$threePmData[9,2] = 0
$threePmData[9,3] = 0
$threePmData[9,4] = 2
$threePmData[9,5] = 0
$threePmData[9,6] = 1
$threePmData[9,7] = 2
$threePmData[9,8] = 0
$threePmData[9,9] = 0
$threePmData[9,10] = 2
$threePmData[9,11] = 1
$threePmData[9,12] = 0
$threePmData[9,13] = 0
$threePmData[9,14] = 1
$threePmData[9,15] = 0
$threePmData[9,16] = 0
$threePmData[10,0] = "2025-11-10"
$threePmData[10,1] = "MIN"
$threePmData[10,2] = 0
$threePmData[10,3] = 0
$threePmData[10,4] = 3
$threePmData[10,5] = 0
$threePmData[10,6] = 0
$threePmData[10,7] = 3
$threePmData[10,8] = 2
$threePmData[10,9] = 0
$threePmData[10,10] = 0
$threePmData[10,11] = 1
$threePmData[10,12] = 4
$threePmData[10,13] = 0
$threePmData[10,14] = 1
$threePmData[10,15] = 1
$threePmData[10,16] = 2
$ws3pm.Range("A1:A11").NumberFormat = "@"
$ws3pm.Range("A1:Q11").Value = $threePmData
$ws3pm.Range("A1:Q1").Font.Bold = $true
$ws3pm.Range("A1:Q1").HorizontalAlignment = -4108
$ws3pm.Range("A1:Q1").VerticalAlignment = -4160
$ws3pm.Range("A1:Q1").Borders.LineStyle = 1

# --- Insert 'Avg Rebounds' sheet right after 'Avg Assists' ---
$afterAvgAssists = $wb.Worksheets.Item("Avg Assists")
$wsAvgRebounds = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterAvgAssists)
$wsAvgRebounds.Name = "Avg Rebounds"

$avgReboundsData = New-Object 'object[,]' 16,2
$avgReboundsData[0,0] = "Player"
$avgReboundsData[0,1] = "Avg Rebounds"
$avgReboundsData[1,0] = "Walker Kessler"
$avgReboundsData[1,1] = 10.8
$avgReboundsData[2,0] = "Jusuf Nurkić"
$avgReboundsData[2,1] = 8.9
$avgReboundsData[3,0] = "Lauri Markkanen"
$avgReboundsData[3,1] = 6.3
$avgReboundsData[4,0] = "Kyle Filipowski"
$avgReboundsData[4,1] = 4.8
$avgReboundsData[5,0] = "Taylor Hendricks"
$avgReboundsData[5,1] = 4.222222222222222
$avgReboundsData[6,0] = "Keyonte George"
$avgReboundsData[6,1] = 3.8
$avgReboundsData[7,0] = "Kevin Love"
$avgReboundsData[7,1] = 3.8
$avgReboundsData[8,0] = "Ace Bailey"
$avgReboundsData[8,1] = 3.7
$avgReboundsData[9,0] = "Svi Mykhailiuk"
$avgReboundsData[9,1] = 2.8
$avgReboundsData[10,0] = "Walter Clayton Jr."
$avgReboundsData[10,1] = 2.6
$avgReboundsData[11,0] = "Kyle Anderson"
$avgReboundsData[11,1] = 2.333333333333333
$avgReboundsData[12,0] = "Brice Sensabaugh"
$avgReboundsData[12,1] = 2.222222222222222
$avgReboundsData[13,0] = "Isaiah Collier"
$avgReboundsData[13,1] = 1.5
$avgReboundsData[14,0] = "Elijah Harkless"
$avgReboundsData[14,1] = 1.2
$avgReboundsData[15,0] = "Cody Williams"
$avgReboundsData[15,1] = 1
$wsAvgRebounds.Range("A1:B16").Value = $avgReboundsData
$wsAvgRebounds.Range("A1:B1").Font.Bold = $true
$wsAvgRebounds.Range("A1:B1").HorizontalAlignment = -4108
$wsAvgRebounds.Range("A1:B1").VerticalAlignment = -4160
$wsAvgRebounds.Range("A1:B1").Borders.LineStyle = 1

# --- Insert 'Avg 3PM' sheet right after 'Avg Rebounds' ---
$wsAvg3pm = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsAvgRebounds)
$wsAvg3pm.Name = "Avg 3PM"

$avg3pmData = New-Object 'object[,]' 16,2
$avg3pmData[0,0] = "Player"
$avg3pmData[0,1] = "Avg 3PM"
$avg3pmData[1,0] = "Lauri Markkanen"
$avg3pmData[1,1] = 3.2
$avg3pmData[2,0] = "Svi Mykhailiuk"
$avg3pmData[2,1] = 2.1
$avg3pmData[3,0] = "Keyonte George"
$avg3pmData[3,1] = 1.7
$avg3pmData[4,0] = "Walker Kessler"
$avg3pmData[4,1] = 1.2
$avg3pmData[5,0] = "Brice Sensabaugh"
$avg3pmData[5,1] = 1.111111111111111
$avg3pmData[6,0] = "Walter Clayton Jr."
$avg3pmData[6,1] = 0.9
$avg3pmData[7,0] = "Kyle Filipowski"
$avg3pmData[7,1] = 0.9
$avg3pmData[8,0] = "Kevin Love"
$avg3pmData[8,1] = 0.6
$avg3pmData[9,0] = "Isaiah Collier"
$avg3pmData[9,1] = 0.5
$avg3pmData[10,0] = "Taylor Hendricks"
$avg3pmData[10,1] = 0.4444444444444444
$avg3pmData[11,0] = "Cody Williams"
$avg3pmData[11,1] = 0.4
$avg3pmData[12,0] = "Elijah Harkless"
$avg3pmData[12,1] = 0.4
$avg3pmData[13,0] = "Ace Bailey"
$avg3pmData[13,1] = 0.4
$avg3pmData[14,0] = "Jusuf Nurkić"
$avg3pmData[14,1] = 0.3
$avg3pmData[15,0] = "Kyle Anderson"
$avg3pmData[15,1] = 0
$wsAvg3pm.Range("A1:B16").Value = $avg3pmData
$wsAvg3pm.Range("A1:B1").Font.Bold = $true
$wsAvg3pm.Range("A1:B1").HorizontalAlignment = -4108
$wsAvg3pm.Range("A1:B1").VerticalAlignment = -4160
$wsAvg3pm.Range("A1:B1").Borders.LineStyle = 1

Write-Output ($wb.Worksheets | ForEach-Object { $_.Name })
